$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gained a new weekly price observation. Insert a new row at
# position 238 (shifting the existing rows 238-246 down to 239-247) and
# populate it with the new record's data.
$ws.Rows(238).Insert()

$ws.Cells.Item(238, 1).Value = 1
$ws.Cells.Item(238, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(238, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(238, 4).Value = 44706
$ws.Cells.Item(238, 5).Value = 15
$ws.Cells.Item(238, 6).Value = "Fruta"
$ws.Cells.Item(238, 7).Value = 100108
$ws.Cells.Item(238, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(238, 9).Value = 100108006
$ws.Cells.Item(238, 10).Value = "Plátano"
$ws.Cells.Item(238, 11).Value = "Sin especificar"
$ws.Cells.Item(238, 12).Value = "Pintón"
$ws.Cells.Item(238, 13).Value = 120
$ws.Cells.Item(238, 14).Value = 14000
$ws.Cells.Item(238, 15).Value = 15000
$ws.Cells.Item(238, 16).Value = 14500
$ws.Cells.Item(238, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(238, 18).Value = "Ecuador"
$ws.Cells.Item(238, 19).Value = 725
$ws.Cells.Item(238, 20).Value = 20
